$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shortage")

$ws.Range("C2").Value = [double]"3.072906383865932e-05"
$ws.Range("D2").Value = [double]"0.0002217231224067183"

$ws.Range("C3").Value = [double]"0.310481901381304"
$ws.Range("D3").Value = [double]"4.134853793981501e-05"

$ws.Range("C4").Value = [double]"5.633554804655667e-05"
$ws.Range("D4").Value = [double]"1.439655649403873e-08"
